$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Insert a new row at the top to make room for the header row
$ws.Rows.Item(1).Insert()

# New header row: Name, Farmer, Ranger, Knight, King
$ws.Cells.Item(1, 1).Value = "Name"
$ws.Cells.Item(1, 2).Value = "Farmer"
$ws.Cells.Item(1, 3).Value = "Ranger"
$ws.Cells.Item(1, 4).Value = "Knight"
$ws.Cells.Item(1, 5).Value = "King"

$ws.Range("E2").Select()
